$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 38, shifting existing rows 38:50 down to 39:51
$ws.Rows.Item(38).Insert()

# Fill in the new row 38 with data
$ws.Cells.Item(38, 1).Value = 8
$ws.Cells.Item(38, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44754
$ws.Cells.Item(38, 5).Value = 4
$ws.Cells.Item(38, 6).Value = 100114007
$ws.Cells.Item(38, 7).Value = "Jengibre"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 480
$ws.Cells.Item(38, 11).Value = 15000
$ws.Cells.Item(38, 12).Value = 16000
$ws.Cells.Item(38, 13).Value = 15500
$ws.Cells.Item(38, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(38, 15).Value = "Perú"
$ws.Cells.Item(38, 16).Value = 1192
$ws.Cells.Item(38, 17).Value = 13
$ws.Cells.Item(38, 18).Value = "Hortaliza"
